# Generate Report for Archive
#
# 1. Update the localization status from "Ready for handoff" to "In Translation"
#    everywhere it appears (Overview, zh-cn and de-de sheets).
# 2. Narrow the "Status" / per-locale status columns to match the new
#    (shorter) text - the Overview sheet's zh-cn/de-de columns (E, F) and the
#    "Status" column (C) on the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rows = $used.Rows.Count
    $cols = $used.Columns.Count
    for ($r = 1; $r -le $rows; $r++) {
        for ($c = 1; $c -le $cols; $c++) {
            $cell = $used.Cells.Item($r, $c)
            if ("Ready for handoff" -eq $cell.Value2) {
                $cell.Value = "In Translation"
            }
        }
    }
}

$newStatusColumnWidth = 12.5

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = $newStatusColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newStatusColumnWidth

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = $newStatusColumnWidth

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = $newStatusColumnWidth
